# MassWateR ParameterMapping.xlsx update
# - Swap "Cyanobacteria (probe)"/"Cyanobacteria (lab)" rows so the plain
#   "Cyanobacteria" simple-parameter (lab) sits on row 37 and the probe
#   variant sits on row 38 (matches WQX parameter + units for each).
# - Fix a typo: "Phcyocyanin (probe)" -> "Phycocyanin (probe)" (row 40, WQX
#   Parameter column).
# - Add a new reviewer comment on D40 explaining why mg/m3 was left out of
#   the units list.
# - Refresh the view state (frozen-pane anchor + active selection) to match
#   the latest authoring session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: now the plain "Cyanobacteria" (lab) entry ---
$ws.Range("B37").Value = "Cyanobacteria"
$ws.Range("C37").Value = "Algae, blue-green (phylum cyanophyta) density"
$ws.Range("D37").Value = "mg/l, ug/l, umol/l, ppm"

# --- Row 38: now the "Cyanobacteria (probe)" entry ---
$ws.Range("B38").Value = "Cyanobacteria (probe)"
$ws.Range("C38").Value = "Chlorophyll a (probe) concentration, Cyanobacteria (bluegreen)"
$ws.Range("D38").Value = "mg/l, ug/l, umol/l, ppm, RFU"

# --- Row 40: fix "Phcyocyanin (probe)" typo in the WQX Parameter column ---
$ws.Range("C40").Value = "Phycocyanin (probe)"

# --- New comment on D40 ---
$comment = $ws.Range("D40").AddComment("Benjamen Wetherill:" + [char]10 + "Did not include mg/m3 because it is the same as ug/l.")

# --- Refresh frozen-pane anchor / active selection to match latest session ---
$ws.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E11").Select()

Write-Host "edit complete"
